$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.82
$ws.Range("G2").Value = 1.85
$ws.Range("H2").Value = 6.8
$ws.Range("I2").Value = 7.6
$ws.Range("J2").Value = 3.15
$ws.Range("K2").Value = 3.3
$ws.Range("L2").Value = 2.6
$ws.Range("M2").Value = 1.18
$ws.Range("N2").Value = 2.14
$ws.Range("O2").Value = 1.82
$ws.Range("P2").Value = 1.34
$ws.Range("Q2").Value = 3.75
$ws.Range("R2").Value = 1.1
$ws.Range("S2").Value = 9.800000000000001
$ws.Range("T2").Value = 3
$ws.Range("U2").Value = 1.45
$ws.Range("V2").Value = 1.16
$ws.Range("W2").Value = 2.12
$ws.Range("X2").Value = 6.4
$ws.Range("Y2").Value = 14
$ws.Range("Z2").Value = 1000
$ws.Range("AB2").Value = 4.9
$ws.Range("AC2").Value = 8.800000000000001
$ws.Range("AD2").Value = 990
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 8.4
$ws.Range("AG2").Value = 15.5
$ws.Range("AH2").Value = 990
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 22
$ws.Range("AK2").Value = 1000
$ws.Range("AL2").Value = 1000
$ws.Range("AN2").Value = 1000
$ws.Range("F3").Value = 1.96
$ws.Range("G3").Value = 1.98
$ws.Range("I3").Value = 4.4
$ws.Range("J3").Value = 3.85
$ws.Range("K3").Value = 4
$ws.Range("N3").Value = 4.2
$ws.Range("O3").Value = 1.29
$ws.Range("P3").Value = 2.04
$ws.Range("Q3").Value = 1.91
$ws.Range("R3").Value = 1.41
$ws.Range("T3").Value = 1.79
$ws.Range("U3").Value = 2.18
$ws.Range("V3").Value = 1.29
$ws.Range("W3").Value = 2.02
$ws.Range("X3").Value = 17.5
$ws.Range("Y3").Value = 16.5
$ws.Range("Z3").Value = 32
$ws.Range("AA3").Value = 140
$ws.Range("AE3").Value = 55
$ws.Range("AH3").Value = 25
$ws.Range("AI3").Value = 85
$ws.Range("AK3").Value = 26
$ws.Range("AM3").Value = 130
$ws.Range("AN3").Value = 15.5
$ws.Range("AO3").Value = 1000
$ws.Range("F4").Value = 1.46
$ws.Range("G4").Value = 1.51
$ws.Range("H4").Value = 6.8
$ws.Range("I4").Value = 7.6
$ws.Range("J4").Value = 5.1
$ws.Range("K4").Value = 5.5
$ws.Range("L4").Value = 1.25
$ws.Range("N4").Value = 6.8
$ws.Range("O4").Value = 1.15
$ws.Range("P4").Value = 2.92
$ws.Range("Q4").Value = 1.47
$ws.Range("T4").Value = 1.7
$ws.Range("U4").Value = 2.28
$ws.Range("V4").Value = 1.15
$ws.Range("W4").Value = 2.96
$ws.Range("X4").Value = 40
$ws.Range("Y4").Value = 980
$ws.Range("Z4").Value = 980
$ws.Range("AA4").Value = 210
$ws.Range("AB4").Value = 23
$ws.Range("AC4").Value = 42
$ws.Range("AD4").Value = 26
$ws.Range("AE4").Value = 260
$ws.Range("AF4").Value = 11.5
$ws.Range("AH4").Value = 980
$ws.Range("AI4").Value = 70
$ws.Range("AJ4").Value = 21
$ws.Range("AK4").Value = 13.5
$ws.Range("AM4").Value = 95
$ws.Range("AN4").Value = 7.4
$ws.Range("AO4").Value = 1000
$ws.Range("F5").Value = 1.53
$ws.Range("N5").Value = 2.34
$ws.Range("O5").Value = 1.73
$ws.Range("P5").Value = 1.44
$ws.Range("R5").Value = 1.13
$ws.Range("T5").Value = 1.04
$ws.Range("U5").Value = 1.43
$ws.Range("V5").Value = 1.09
$ws.Range("X5").Value = 9
$ws.Range("Y5").Value = 21
$ws.Range("Z5").Value = 1000
$ws.Range("AB5").Value = 4.7
$ws.Range("AD5").Value = 60
$ws.Range("AI5").Value = 620
$ws.Range("AK5").Value = 27
$ws.Range("AL5").Value = 120
$ws.Range("AM5").Value = 870
$ws.Range("F6").Value = 1.71
$ws.Range("G6").Value = 1.79
$ws.Range("H6").Value = 5.2
$ws.Range("I6").Value = 6.8
$ws.Range("J6").Value = 3.8
$ws.Range("K6").Value = 4.3
$ws.Range("O6").Value = 1.35
$ws.Range("P6").Value = 1.84
$ws.Range("R6").Value = 1.3
$ws.Range("S6").Value = 3.75
$ws.Range("T6").Value = 1.05
$ws.Range("V6").Value = 1.19
$ws.Range("W6").Value = 2.26
$ws.Range("X6").Value = 1000
$ws.Range("Y6").Value = 1000
$ws.Range("Z6").Value = 1000
$ws.Range("AB6").Value = 29
$ws.Range("AC6").Value = 42
$ws.Range("AD6").Value = 1000
$ws.Range("AE6").Value = 1000
$ws.Range("AF6").Value = 980
$ws.Range("AG6").Value = 40
$ws.Range("AH6").Value = 1000
$ws.Range("AJ6").Value = 1000
$ws.Range("AK6").Value = 1000
$ws.Range("AL6").Value = 1000
$ws.Range("AN6").Value = 85
$ws.Range("F7").Value = 1.67
$ws.Range("G7").Value = 1.71
$ws.Range("H7").Value = 5.8
$ws.Range("I7").Value = 7
$ws.Range("J7").Value = 3.9
$ws.Range("K7").Value = 4.2
$ws.Range("L7").Value = 1.42
$ws.Range("N7").Value = 3.4
$ws.Range("O7").Value = 1.35
$ws.Range("Q7").Value = 2.06
$ws.Range("R7").Value = 1.31
$ws.Range("T7").Value = 1.05
$ws.Range("V7").Value = 1.17
$ws.Range("W7").Value = 2.4
$ws.Range("Y7").Value = 19.5
$ws.Range("Z7").Value = 980
$ws.Range("AA7").Value = 900
$ws.Range("AB7").Value = 7.6
$ws.Range("AC7").Value = 42
$ws.Range("AD7").Value = 25
$ws.Range("AF7").Value = 9
$ws.Range("AH7").Value = 25
$ws.Range("AI7").Value = 190
$ws.Range("AJ7").Value = 16.5
$ws.Range("AK7").Value = 19.5
$ws.Range("AL7").Value = 70
$ws.Range("AN7").Value = 85
$ws.Range("AO7").Value = 1000
$ws.Range("F8").Value = 2.2
$ws.Range("H8").Value = 3.55
$ws.Range("I8").Value = 3.95
$ws.Range("N8").Value = 3.3
$ws.Range("O8").Value = 1.4
$ws.Range("P8").Value = 1.76
$ws.Range("Q8").Value = 2.14
$ws.Range("R8").Value = 1.3
$ws.Range("S8").Value = 4.1
$ws.Range("T8").Value = 1.05
$ws.Range("U8").Value = 1.04
$ws.Range("V8").Value = 1.34
$ws.Range("X8").Value = 14.5
$ws.Range("Y8").Value = 1000
$ws.Range("Z8").Value = 980
$ws.Range("AB8").Value = 1000
$ws.Range("AC8").Value = 42
$ws.Range("AD8").Value = 980
$ws.Range("AE8").Value = 980
$ws.Range("AF8").Value = 1000
$ws.Range("AG8").Value = 40
$ws.Range("AH8").Value = 60
$ws.Range("AI8").Value = 330
$ws.Range("AJ8").Value = 980
$ws.Range("AK8").Value = 980
$ws.Range("AL8").Value = 980
$ws.Range("AM8").Value = 1000
$ws.Range("AN8").Value = 600
$ws.Range("AO8").Value = 1000
$ws.Range("F9").Value = 2.52
$ws.Range("G9").Value = 2.54
$ws.Range("L9").Value = 1.52
$ws.Range("N9").Value = 3.1
$ws.Range("O9").Value = 1.44
$ws.Range("P9").Value = 1.69
$ws.Range("Q9").Value = 2.38
$ws.Range("R9").Value = 1.25
$ws.Range("S9").Value = 4.5
$ws.Range("T9").Value = 1.96
$ws.Range("U9").Value = 1.96
$ws.Range("V9").Value = 1.44
$ws.Range("W9").Value = 1.64
$ws.Range("X9").Value = 11.5
$ws.Range("Y9").Value = 12
$ws.Range("AA9").Value = 500
$ws.Range("AC9").Value = 7
$ws.Range("AD9").Value = 26
$ws.Range("AE9").Value = 44
$ws.Range("AG9").Value = 12
$ws.Range("AH9").Value = 20
$ws.Range("AI9").Value = 65
$ws.Range("AL9").Value = 55
$ws.Range("AM9").Value = 580
$ws.Range("AN9").Value = 990
$ws.Range("AO9").Value = 970
$ws.Range("F10").Value = 1.98
$ws.Range("G10").Value = 2.02
$ws.Range("K10").Value = 3.65
$ws.Range("O10").Value = 1.39
$ws.Range("P10").Value = 1.81
$ws.Range("R10").Value = 1.3
$ws.Range("S10").Value = 3.95
$ws.Range("T10").Value = 1.89
$ws.Range("W10").Value = 1.99
$ws.Range("X10").Value = 12
$ws.Range("Y10").Value = 15
$ws.Range("Z10").Value = 85
$ws.Range("AB10").Value = 10.5
$ws.Range("AC10").Value = 11.5
$ws.Range("AD10").Value = 18
$ws.Range("AJ10").Value = 22
$ws.Range("AK10").Value = 22
$ws.Range("AL10").Value = 42
$ws.Range("AN10").Value = 17.5
$ws.Range("AO10").Value = 80
